$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 3
$ws.Range("E28").Value = 11
$ws.Range("E33").Value = 28
$ws.Range("F33").Value = 8
$ws.Range("H33").Value = 8
$ws.Range("E37").Value = 37
$ws.Range("E61").Value = 22
$ws.Range("E67").Value = 33
$ws.Range("E77").Value = 43
$ws.Range("E79").Value = 25
